$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 32.3
$ws.Range("B24").Value = 8.4
$ws.Range("C24").Value = 26.2
$ws.Range("D24").Value = 6.4
$ws.Range("E24").Value = 24.2
$ws.Range("F24").Value = 12.6
$ws.Range("G24").Value = 12.1
$ws.Range("H24").Value = 10.9
$ws.Range("I24").Value = 7.4
$ws.Range("J24").Value = 15.9
$ws.Range("K24").Value = 44046.91666666666
$ws.Range("K24").NumberFormat = $ws.Range("K23").NumberFormat
